$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: for numeric-looking text values (e.g. "596.47"), Excel
# auto-converts Range.Value to a Double. Forcing NumberFormat to text ("@")
# before the assignment keeps it a string; resetting Style to "Normal"
# afterwards drops the explicit text number-format so the cell keeps the
# workbook's original (unstyled) look, matching the source data files produced
# by the scraping script.

$ws.Range('D2').Value = '67.571.90'
$ws.Range('E2').Value = '  -2.38%  '

$ws.Range('D3').Value = '3.797.77'
$ws.Range('E3').Value = '  +1.58%  '

$ws.Range('E4').Value = '  +0.32%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '596.47'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -2.73%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '174.25'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.10%  '

$ws.Range('D7').Value = '3.794.55'
$ws.Range('E7').Value = '  +1.56%  '

$ws.Range('E8').Value = '  +0.14%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.518'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -1.35%  '

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.158'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -4.21%  '

$ws.Range('E11').Value = '  -5.74%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.461'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -3.81%  '

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '37.94'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -4.71%  '

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.0000244'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -3.43%  '

$ws.Range('D15').Value = '4.442.92'
$ws.Range('E15').Value = '  +1.88%  '

$ws.Range('D16').Value = '3.833.63'
$ws.Range('E16').Value = '  +2.61%  '

$ws.Range('D17').Value = '67.725.10'
$ws.Range('E17').Value = '  -2.25%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.115'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -4.69%  '

$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.14'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -3.57%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '16.26'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -0.24%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '487.04'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -2.26%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '9.11'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.60%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.726'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +1.04%  '

$ws.Range('E24').Value = '  +11.60%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '83.56'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -2.32%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -7.68%  '

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '12.19'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -5.22%  '

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '10.21'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -6.13%  '

$ws.Range('E29').Value = '  +0.25%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.95'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.59%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -2.26%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '32.94'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +8.79%  '

$ws.Range('E33').Value = '  -3.54%  '

$ws.Range('E34').Value = '  -3.42%  '

$ws.Range('E35').Value = '  +0.47%  '

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -3.54%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.136'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -0.71%  '

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '5.74'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -5.75%  '

$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.328'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -5.29%  '

$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '454.91'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +2.82%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '49.11'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -1.15%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -2.84%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.85'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -7.00%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '8.29'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -2.76%  '

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '41.38'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -6.56%  '

$ws.Range('D46').Value = '2.825.87'
$ws.Range('E46').Value = '  -4.02%  '

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '140.63'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +1.88%  '

$ws.Range('E48').Value = '  +0.01%  '

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.0349'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -2.37%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '25.97'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -3.82%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '23.52'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +9.85%  '
